$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, matching the style/formatting of the
# neighboring header cell (G1: bold, bordered, centered/top-aligned)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add the corresponding value for the new column
$ws.Range("H2").Value = 1
